# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4).
#
#   F2: 10022 -> 10061
#   F3:   222 ->   223
#   F4:    46 ->    47
#   F5:   606 ->   611

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 10061
    $ws.Range("F3").Value = 223
    $ws.Range("F4").Value = 47
    $ws.Range("F5").Value = 611
}
